$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the large "H" column numeric values (series data)
$ws.Range("H9").Value  = 78974589789
$ws.Range("H11").Value = 43534545345345
$ws.Range("H12").Value = 3454545345345
$ws.Range("H14").Value = 345345454534534
$ws.Range("H15").Value = 45345345345345
$ws.Range("H16").Value = 343453454534545
$ws.Range("H18").Value = 34545345345454
$ws.Range("H19").Value = 3434543345
$ws.Range("H20").Value = 3454534534543529984

# Move the active selection / view back to the top of the sheet (H4),
# which also clears the previous scrolled-down topLeftCell="G19" state.
$ws.Range("H4").Select()
